$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"3.51660537190763E-06"
$ws.Range("E2").Value = [double]"3.51660537190763E-06"

$ws.Range("D3").Value = [double]"5.647895547456047E-23"
$ws.Range("E3").Value = [double]"5.647895547456047E-23"

$ws.Range("D4").Value = [double]"1.307617941001401E-11"
$ws.Range("E4").Value = [double]"1.307617941001401E-11"

$ws.Range("D5").Value = [double]"8.341119027495172E-17"
$ws.Range("E5").Value = [double]"8.341119027495172E-17"

$ws.Range("D6").Value = [double]"0.9999999980978174"
$ws.Range("E6").Value = [double]"0.9999999980978174"

$ws.Range("D8").Value = [double]"5.904795123004318E-05"
$ws.Range("E8").Value = [double]"0.99994095204877"

$ws.Range("D9").Value = [double]"0.122122514494204"
$ws.Range("E9").Value = [double]"0.877877485505796"

$ws.Range("D11").Value = [double]"0.9999999999999791"
$ws.Range("E11").Value = [double]"2.087219286295294E-14"
$ws.Range("F11").Value = [double]"3.192015886306763"
